$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '25.804.20'
$ws.Range("E2").Value = '  +0.45%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.749.86'
$ws.Range("E3").Value = '  +0.21%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '236.10'
$ws.Range("E5").Value = '  -0.29%  '
$ws.Range("E6").Value = '  -0.02%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5089'
$ws.Range("E7").Value = '  +3.70%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '40.77'
$ws.Range("E8").Value = '  -1.95%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.2673'
$ws.Range("E9").Value = '  +7.52%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.06195'
$ws.Range("E10").Value = '  +3.85%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.751.99'
$ws.Range("E11").Value = '  +0.42%  '
$ws.Range("E12").Value = '  +2.21%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '15.43'
$ws.Range("E13").Value = '  +4.38%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.6279'
$ws.Range("E14").Value = '  +11.99%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.475'
$ws.Range("E15").Value = '  +0.38%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '77.74'
$ws.Range("E16").Value = '  +0.72%  '
$ws.Range("E18").Value = '  -0.02%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '25.812.24'
$ws.Range("E19").Value = '  +0.33%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.64'
$ws.Range("E20").Value = '  +1.66%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.000006661'
$ws.Range("E21").Value = '  +1.84%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.978.27'
$ws.Range("E22").Value = '  +0.67%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.057'
$ws.Range("E23").Value = '  +2.31%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '8.259'
$ws.Range("E24").Value = '  +5.22%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '5.132'
$ws.Range("E25").Value = '  +2.50%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '136.75'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.456'
$ws.Range("E27").Value = '  -1.73%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '15.11'
$ws.Range("E28").Value = '  +3.32%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.744'
$ws.Range("E29").Value = '  -3.54%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '102.50'
$ws.Range("E30").Value = '  +0.64%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08188'
$ws.Range("E31").Value = '  +2.15%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.701'
$ws.Range("E32").Value = '  -1.11%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.403'
$ws.Range("E33").Value = '  +3.09%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.04419'
$ws.Range("E34").Value = '  +0.77%  '
$ws.Range("B35").Value = 'HuobiToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.651'
$ws.Range("E35").Value = '  +2.53%  '
$ws.Range("B36").Value = 'ARBITRUM'
$ws.Range("C36").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.9972'
$ws.Range("E36").Value = '  +1.35%  '
$ws.Range("B37").Value = 'ImmutableX'
$ws.Range("C37").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.5989'
$ws.Range("E37").Value = '  -0.81%  '
$ws.Range("B38").Value = 'MXToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.645'
$ws.Range("E38").Value = '  -1.29%  '
$ws.Range("B39").Value = 'VeChain'
$ws.Range("C39").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01565'
$ws.Range("E39").Value = '  +4.55%  '
$ws.Range("B40").Value = 'RenderToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.943'
$ws.Range("E40").Value = '  -3.24%  '
$ws.Range("B41").Value = 'PaxDollar'
$ws.Range("C41").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.001'
$ws.Range("E41").Value = '  +0.00%  '
$ws.Range("B42").Value = 'Quant'
$ws.Range("C42").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '101.52'
$ws.Range("E42").Value = '  -1.84%  '
$ws.Range("B43").Value = 'TrustWalletToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.7510'
$ws.Range("E43").Value = '  -0.71%  '
$ws.Range("B44").Value = 'TheSandbox'
$ws.Range("C44").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.3823'
$ws.Range("E44").Value = '  +3.32%  '
$ws.Range("B45").Value = 'FraxShare'
$ws.Range("C45").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '4.891'
$ws.Range("E45").Value = '  -4.79%  '
$ws.Range("B46").Value = 'Cronos'
$ws.Range("C46").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.05501'
$ws.Range("E46").Value = '  +7.65%  '
$ws.Range("B47").Value = 'Algorand'
$ws.Range("C47").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.1097'
$ws.Range("E47").Value = '  +2.96%  '
$ws.Range("B48").Value = 'Aptos'
$ws.Range("C48").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '5.938'
$ws.Range("E48").Value = '  +1.38%  '
$ws.Range("B49").Value = 'Elrond'
$ws.Range("C49").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '30.06'
$ws.Range("E49").Value = '  +0.11%  '
$ws.Range("B50").Value = 'Aave'
$ws.Range("C50").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '52.71'
$ws.Range("E50").Value = '  +0.53%  '
$ws.Range("B51").Value = 'USDD'
$ws.Range("C51").Value = 'https://coinranking.com/coin/z2PZIKQL7+usdd-usdd'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.004'
$ws.Range("E51").Value = '  +0.49%  '
